$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: Camberwell -> Caulfield (Metro Train entry)
$ws.Range("A7").Value = "Caulfield"
$ws.Range("B7").Value = "Metro Train - Frankston line"
$ws.Range("C7").Value = "30/12/20 4:30pm-5:00pm"
$ws.Range("D7").Value = "Case caught train from Caulfield to Cheltenham"

# Row 8: Caulfield -> Cheltenham (Two Bob Snob entry)
$ws.Range("A8").Value = "Cheltenham"
$ws.Range("B8").Value = "Two Bob Snob, 256 Charman Road"
$ws.Range("C8").Value = "22/12/2020 1pm - 2pm"
$ws.Range("D8").Value = "Case attended Venue"

# Row 9: Cheltenham -> Clayton South (Metro Train Pakenham line entry)
$ws.Range("A9").Value = "Clayton South"
$ws.Range("B9").Value = "Metro Train - Pakenham line"
$ws.Range("C9").Value = "31/12/20 9:00pm-9:30pm"
$ws.Range("D9").Value = "Case caught train from Westall Station to Flinders St."

# Row 27: update exposure period
$ws.Range("C27").Value = "01/01/21 4:30am-5:00am"

# Row 28: update exposure period
$ws.Range("C28").Value = "01/01/2021 2:00am - 2:30am"
